$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.01139766666666667
$ws.Range("H2").Value = 0.034193
$ws.Range("I2").Value = 0.1481191086775714
$ws.Range("J2").Value = 0.1481191086775714
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.771611
$ws.Range("N2").Value = 5.314833
$ws.Range("O2").Value = 0.03672985187529028
$ws.Range("P2").Value = 0.03672985187529029
$ws.Range("Q2").Value = 0.020192231641
$ws.Range("R2").Value = 0.181730084769
$ws.Range("S2").Value = 0.00544039292162722
$ws.Range("T2").Value = 0.005440392921627222
$ws.Range("G3").Value = 0.01139766666666667
$ws.Range("H3").Value = 0.034193
$ws.Range("I3").Value = 0.1481191086775714
$ws.Range("J3").Value = 0.1481191086775714
$ws.Range("O3").Value = 0.8428101954878733
$ws.Range("P3").Value = 0.8428101954878736
$ws.Range("Q3").Value = 0.4633348033765556
$ws.Range("R3").Value = 4.170013230389
$ws.Range("S3").Value = 0.1248362949400335
$ws.Range("T3").Value = 0.1248362949400335
$ws.Range("G4").Value = 0.01139766666666667
$ws.Range("H4").Value = 0.034193
$ws.Range("I4").Value = 0.1481191086775714
$ws.Range("J4").Value = 0.1481191086775714
$ws.Range("M4").Value = 0.2535886666666667
$ws.Range("N4").Value = 0.7607660000000001
$ws.Range("O4").Value = 0.005257516556354092
$ws.Range("P4").Value = 0.005257516556354093
$ws.Range("Q4").Value = 0.002890319093111111
$ws.Range("R4").Value = 0.026012871838
$ws.Range("S4").Value = 0.0007787386661847425
$ws.Range("T4").Value = 0.0007787386661847429
$ws.Range("G5").Value = 0.01139766666666667
$ws.Range("H5").Value = 0.034193
$ws.Range("I5").Value = 0.1481191086775714
$ws.Range("J5").Value = 0.1481191086775714
$ws.Range("M5").Value = 5.336666
$ws.Range("N5").Value = 16.009998
$ws.Range("O5").Value = 0.1106422073964871
$ws.Range("P5").Value = 0.1106422073964871
$ws.Range("Q5").Value = 0.06082554017933334
$ws.Range("R5").Value = 0.547429861614
$ws.Range("S5").Value = 0.01638822514168666
$ws.Range("T5").Value = 0.01638822514168667
$ws.Range("G6").Value = 0.01139766666666667
$ws.Range("H6").Value = 0.034193
$ws.Range("I6").Value = 0.1481191086775714
$ws.Range("J6").Value = 0.1481191086775714
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.219956
$ws.Range("N6").Value = 0.659868
$ws.Range("O6").Value = 0.004560228683995159
$ws.Range("P6").Value = 0.00456022868399516
$ws.Range("Q6").Value = 0.002506985169333334
$ws.Range("R6").Value = 0.022562866524
$ws.Range("S6").Value = 0.0006754570080392574
$ws.Range("T6").Value = 0.0006754570080392575
$ws.Range("G7").Value = 0.015206
$ws.Range("H7").Value = 0.045618
$ws.Range("I7").Value = 0.197610548932631
$ws.Range("J7").Value = 0.197610548932631
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.771611
$ws.Range("N7").Value = 5.314833
$ws.Range("O7").Value = 0.03672985187529028
$ws.Range("P7").Value = 0.03672985187529029
$ws.Range("Q7").Value = 0.026939116866
$ws.Range("R7").Value = 0.242452051794
$ws.Range("S7").Value = 0.007258206191290337
$ws.Range("T7").Value = 0.007258206191290341
$ws.Range("G8").Value = 0.015206
$ws.Range("H8").Value = 0.045618
$ws.Range("I8").Value = 0.197610548932631
$ws.Range("J8").Value = 0.197610548932631
$ws.Range("O8").Value = 0.8428101954878733
$ws.Range("P8").Value = 0.8428101954878736
$ws.Range("Q8").Value = 0.6181501202126667
$ws.Range("R8").Value = 5.563351081914
$ws.Range("S8").Value = 0.1665481853763767
$ws.Range("T8").Value = 0.1665481853763768
$ws.Range("G9").Value = 0.015206
$ws.Range("H9").Value = 0.045618
$ws.Range("I9").Value = 0.197610548932631
$ws.Range("J9").Value = 0.197610548932631
$ws.Range("M9").Value = 0.2535886666666667
$ws.Range("N9").Value = 0.7607660000000001
$ws.Range("O9").Value = 0.005257516556354092
$ws.Range("P9").Value = 0.005257516556354093
$ws.Range("Q9").Value = 0.003856069265333333
$ws.Range("R9").Value = 0.034704623388
$ws.Range("S9").Value = 0.001038940732723528
$ws.Range("T9").Value = 0.001038940732723528
$ws.Range("G10").Value = 0.015206
$ws.Range("H10").Value = 0.045618
$ws.Range("I10").Value = 0.197610548932631
$ws.Range("J10").Value = 0.197610548932631
$ws.Range("M10").Value = 5.336666
$ws.Range("N10").Value = 16.009998
$ws.Range("O10").Value = 0.1106422073964871
$ws.Range("P10").Value = 0.1106422073964871
$ws.Range("Q10").Value = 0.08114934319599999
$ws.Range("R10").Value = 0.7303440887639999
$ws.Range("S10").Value = 0.02186406733873781
$ws.Range("T10").Value = 0.02186406733873782
$ws.Range("G11").Value = 0.015206
$ws.Range("H11").Value = 0.045618
$ws.Range("I11").Value = 0.197610548932631
$ws.Range("J11").Value = 0.197610548932631
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.219956
$ws.Range("N11").Value = 0.659868
$ws.Range("O11").Value = 0.004560228683995159
$ws.Range("P11").Value = 0.00456022868399516
$ws.Range("Q11").Value = 0.003344650936
$ws.Range("R11").Value = 0.030101858424
$ws.Range("S11").Value = 0.0009011492935026128
$ws.Range("T11").Value = 0.0009011492935026131
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.05034566666666667
$ws.Range("H12").Value = 0.151037
$ws.Range("I12").Value = 0.6542703423897976
$ws.Range("J12").Value = 0.6542703423897976
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.771611
$ws.Range("N12").Value = 5.314833
$ws.Range("O12").Value = 0.03672985187529028
$ws.Range("P12").Value = 0.03672985187529029
$ws.Range("Q12").Value = 0.08919293686900001
$ws.Range("R12").Value = 0.8027364318210001
$ws.Range("S12").Value = 0.02403125276237272
$ws.Range("T12").Value = 0.02403125276237273
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.05034566666666667
$ws.Range("H13").Value = 0.151037
$ws.Range("I13").Value = 0.6542703423897976
$ws.Range("J13").Value = 0.6542703423897976
$ws.Range("O13").Value = 0.8428101954878733
$ws.Range("P13").Value = 0.8428101954878736
$ws.Range("Q13").Value = 2.046638162711222
$ws.Range("R13").Value = 18.419743464401
$ws.Range("S13").Value = 0.5514257151714631
$ws.Range("T13").Value = 0.5514257151714633
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05034566666666667
$ws.Range("H14").Value = 0.151037
$ws.Range("I14").Value = 0.6542703423897976
$ws.Range("J14").Value = 0.6542703423897976
$ws.Range("M14").Value = 0.2535886666666667
$ws.Range("N14").Value = 0.7607660000000001
$ws.Range("O14").Value = 0.005257516556354092
$ws.Range("P14").Value = 0.005257516556354093
$ws.Range("Q14").Value = 0.01276709048244445
$ws.Range("R14").Value = 0.114903814342
$ws.Range("S14").Value = 0.003439837157445821
$ws.Range("T14").Value = 0.003439837157445822
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05034566666666667
$ws.Range("H15").Value = 0.151037
$ws.Range("I15").Value = 0.6542703423897976
$ws.Range("J15").Value = 0.6542703423897976
$ws.Range("M15").Value = 5.336666
$ws.Range("N15").Value = 16.009998
$ws.Range("O15").Value = 0.1106422073964871
$ws.Range("P15").Value = 0.1106422073964871
$ws.Range("Q15").Value = 0.2686780075473333
$ws.Range("R15").Value = 2.418102067926
$ws.Range("S15").Value = 0.0723899149160626
$ws.Range("T15").Value = 0.07238991491606261
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05034566666666667
$ws.Range("H16").Value = 0.151037
$ws.Range("I16").Value = 0.6542703423897976
$ws.Range("J16").Value = 0.6542703423897976
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.219956
$ws.Range("N16").Value = 0.659868
$ws.Range("O16").Value = 0.004560228683995159
$ws.Range("P16").Value = 0.00456022868399516
$ws.Range("Q16").Value = 0.01107383145733334
$ws.Range("R16").Value = 0.099664483116
$ws.Range("S16").Value = 0.00298362238245329
$ws.Range("T16").Value = 0.00298362238245329

Write-Output "Updated 196 cells"
